# Adding in genus/organisms for Verheyen et al 2020
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use an existing cell that already carries the "s=2" style (font with explicit
# black color) so that newly written numeric cells reuse that same style index
# instead of Excel creating a brand new one.
$styleFont = $ws.Range("W274").Font.Color

for ($r = 274; $r -le 369; $r++) {
    $ws.Cells.Item($r, 24).Value = "Ischnura "   # column X - genus
    $ws.Cells.Item($r, 25).Value = "elegans "    # column Y - species

    $zCell  = $ws.Cells.Item($r, 26)             # column Z  - larger_group
    $aaCell = $ws.Cells.Item($r, 27)              # column AA - exp_age
    $abCell = $ws.Cells.Item($r, 28)              # column AB - size

    $zCell.Font.Color = $styleFont
    $aaCell.Font.Color = $styleFont
    $abCell.Font.Color = $styleFont

    $zCell.Value = 1
    $aaCell.Value = 0
    $abCell.Value = 1
}

# Best-effort: restore the selection/scroll state that was recorded after the edit.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 330
$ws.Range("N338").Select()
